# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from 2023-09-01 (45170) to 2023-09-05 (45174), matching the
# diff's change of the underlying date serial value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).Date.AddDays(45174)

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
